$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update R2 row (row 6): swap the 232k note/part for the 226k one ---
$ws.Range("C6").Value = "311-226KFRCT-ND"
$ws.Range("D6").Value = "Using 226k resistor which is a standard size down"

# --- Add new row 10: Fuse Holder ---
# Copy formatting from the row above (row 9) first so the new row gets the
# same per-column style indices (5 for B/C/D/H, 6 for E/F/G/I/J), then
# overwrite with the real values.
$ws.Range("B9:J9").Copy() | Out-Null
$ws.Range("B10:J10").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "F5187-ND"
$ws.Range("E10").Value = "ATO_FUSEHOLDER"
$ws.Range("F10").Value = "FUSEATO_FUSEHOLDER"
$ws.Range("G10").Value = "ATO_BLADE_FUSE"
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = "FUSE"
$ws.Range("J10").Value = "Fuse Holder"

# --- Add new row 11: SPST Relay ---
$ws.Range("B9:J9").Copy() | Out-Null
$ws.Range("B11:J11").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "Z962-ND"
$ws.Range("E11").Value = "G4W-1114P-US-TV8-HP"
$ws.Range("F11").Value = "G4W-1114P-US-TV8-HP"
$ws.Range("G11").Value = "G4W-1114P-US-TV8-HP"
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = "RELAY"
$ws.Range("J11").Value = "SPST Relay"

$excel.CutCopyMode = 0

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 14.28515625
$ws.Columns.Item(3).ColumnWidth = 25.42578125
$ws.Columns.Item(5).ColumnWidth = 21.5703125
$ws.Columns.Item(6).ColumnWidth = 21.7109375
$ws.Columns.Item(7).ColumnWidth = 21.5703125
$ws.Columns.Item(8).ColumnWidth = 8.5703125

# --- AutoFilter on the header row ---
$ws.Range("B2:J2").AutoFilter()
